# fix(publipostage): Correct status name
#
# The commit renames a handful of shared-string values used by the
# "statut_label" (column B) and "statut_name" (column C) columns:
#   "bleu"                                                  -> "noir"
#   "pas de résultat ni de publication"                     -> "pas de résultat postés ni publiés"
#   "résultat et / ou publication posté"                    -> "résultat postés ou publiés"
#   "résultat et / ou publication posté dans les 36 mois"   -> "résultat postés ou publiés dans les 36 mois"
#   "résultat et / ou publication posté dans les 12 mois"   -> "résultat postés ou publiés dans les 12 mois"
#
# Use Range.Replace with LookAt:=xlWhole so only cells whose entire
# content equals the old label are updated (avoids touching unrelated
# cells that merely contain similar substrings, e.g. English trial
# titles containing "post").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlWhole = 1      # Excel.XlLookAt.xlWhole
$xlByRows = 1     # Excel.XlSearchOrder.xlByRows
$xlNext = 1       # Excel.XlSearchDirection.xlNext

$replacements = @(
    @{ Old = "bleu"; New = "noir" },
    @{ Old = "pas de résultat ni de publication"; New = "pas de résultat postés ni publiés" },
    @{ Old = "résultat et / ou publication posté"; New = "résultat postés ou publiés" },
    @{ Old = "résultat et / ou publication posté dans les 36 mois"; New = "résultat postés ou publiés dans les 36 mois" },
    @{ Old = "résultat et / ou publication posté dans les 12 mois"; New = "résultat postés ou publiés dans les 12 mois" }
)

foreach ($r in $replacements) {
    $ws.Cells.Replace(
        $r.Old,
        $r.New,
        $xlWhole,
        $xlByRows,
        $false,
        $false,
        $false,
        $false
    )
}
